$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CBM thickness")

# Capitalize "Mid-", "Center" and "Peripheral" in Cuthbertson & Mandel, 1986 references.
# Cell A7 (Balb/c, 1.5 mo, mid-zone) is edited last to reproduce the shared-string
# append order observed in the target workbook.
$ws.Range("A2").Value = "Cuthbertson & Mandel, 1986 (1.5 mo. CBA mice & Retina, Mid-zone)"
$ws.Range("A3").Value = "Cuthbertson & Mandel, 1986 (4 mo. CBA mice & Retina, Mid-zone)"
$ws.Range("A4").Value = "Cuthbertson & Mandel, 1986 (8 mo. CBA mice & Retina, Mid-zone)"
$ws.Range("A5").Value = "Cuthbertson & Mandel, 1986 (12 mo. CBA mice & Retina, Mid-zone)"
$ws.Range("A6").Value = "Cuthbertson & Mandel, 1986 (20 mo. CBA mice & Retina, Mid-zone)"
$ws.Range("A8").Value = "Cuthbertson & Mandel, 1986 (8 mo. Balb/c mice & Retina, Mid-zone)"
$ws.Range("A9").Value = "Cuthbertson & Mandel, 1986 (20 mo. Balb/c mice & Retina, Mid-zone)"
$ws.Range("A10").Value = "Cuthbertson & Mandel, 1986 (1.5 mo. Balb/c mice & Retina, Center zone)"
$ws.Range("A11").Value = "Cuthbertson & Mandel, 1986 (1.5 mo. Balb/c mice & Retina, Periphery zone)"
$ws.Range("A7").Value = "Cuthbertson & Mandel, 1986 (1.5 mo. Balb/c mice & Retina, Mid-zone)"

# Give I33 an explicit formula (matches B33 * 2 / 100, same cached result as before).
$ws.Range("I33").Formula = "=B33*2/100"

# Move the selection/cursor to A8 (also drops the stale topLeftCell scroll anchor).
$ws.Range("A8").Select()
